# Update "想去人数" (want-to-go count) figures in column F across sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 319
$ws1.Range("F4").Value = 1286
$ws1.Range("F5").Value = 635

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 11

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 319
$ws4.Range("F4").Value = 1286
$ws4.Range("F5").Value = 11
$ws4.Range("F6").Value = 635

$wb.Save()
